# The commit swaps the presentation's theme content: the deck's live/active
# theme part (the one the slide master actually points at) had the "Integral"
# green palette and, after the edit, carries the "Office Theme" blue palette
# (the palette that used to live in the sibling, unreferenced theme part).
#
# Concretely the 12 scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) of the active theme change from the Integral values:
#   dk1 000000  lt1 FFFFFF  dk2 455F51  lt2 E3DED1
#   accent1 99CB38 accent2 63A537 accent3 E6D024 accent4 CC9700
#   accent5 4EB3CF accent6 378DA6 hlink 6B9F25 folHlink B26B02
# to the Office Theme values:
#   dk1 000000  lt1 FFFFFF  dk2 44546A  lt2 E7E6E6
#   accent1 5B9BD5 accent2 ED7D31 accent3 A5A5A5 accent4 FFC000
#   accent5 4472C4 accent6 70AD47 hlink 0563C1 folHlink 954F72
#
# The font scheme and format scheme are already identical between the two
# theme parts, so updating the color scheme reproduces the visible effect of
# the swap. PowerPoint's ThemeColorScheme exposes exactly these 12 slots, in
# this order, as a 1-based collection, and RGB is the standard OLE BGR-packed
# integer (0x00BBGGRR) used throughout the COM color APIs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$officeTheme = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Item($i).RGB = ToOleRgb $officeTheme[$i - 1]
}
